$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1,1).Range.Text = "27÷8="
$t.Cell(1,2).Range.Text = "28÷5="
$t.Cell(1,3).Range.Text = "14÷6="
$t.Cell(1,4).Range.Text = "75÷4="
$t.Cell(1,5).Range.Text = "97÷6="

# Row 5
$t.Cell(5,1).Range.Text = "61÷2="
$t.Cell(5,2).Range.Text = "75÷4="
$t.Cell(5,3).Range.Text = "18÷9="
$t.Cell(5,4).Range.Text = "41÷3="
$t.Cell(5,5).Range.Text = "47÷6="

# Row 9
$t.Cell(9,1).Range.Text = "17÷8="
$t.Cell(9,2).Range.Text = "40÷2="
$t.Cell(9,3).Range.Text = "71÷2="
$t.Cell(9,4).Range.Text = "48÷4="
$t.Cell(9,5).Range.Text = "25÷9="

# Row 13
$t.Cell(13,1).Range.Text = "69÷5="
$t.Cell(13,2).Range.Text = "13÷2="
$t.Cell(13,3).Range.Text = "24÷2="
$t.Cell(13,4).Range.Text = "16÷5="
$t.Cell(13,5).Range.Text = "51÷7="

# Row 17
$t.Cell(17,1).Range.Text = "13÷4="
$t.Cell(17,2).Range.Text = "48÷5="
$t.Cell(17,3).Range.Text = "81÷6="
$t.Cell(17,4).Range.Text = "29÷3="
$t.Cell(17,5).Range.Text = "85÷9="

Write-Host "Done updating 25 cells"
